$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 5
$ws.Range("H5").Value = 296
$ws.Range("I5").Value = 231.5
$ws.Range("J5").Value = 425
$ws.Range("K5").Value = 231.5
$ws.Range("L5").Value = 425
$ws.Range("M5").Value = -116.5
$ws.Range("N5").Value = -655
# Row 6
$ws.Range("H6").Value = 298.125
$ws.Range("I6").Value = 280.16666
$ws.Range("J6").Value = 352
$ws.Range("K6").Value = 840.4999799999999
$ws.Range("L6").Value = 1056
$ws.Range("M6").Value = -728.4999799999999
$ws.Range("N6").Value = -1280
# Row 8
$ws.Range("H8").Value = 58.6
$ws.Range("I8").Value = 64.666664
$ws.Range("J8").Value = 49.5
$ws.Range("K8").Value = 193.999992
$ws.Range("L8").Value = 148.5
$ws.Range("M8").Value = -54.99999199999999
$ws.Range("N8").Value = -426.5
# Row 19
$ws.Range("H19").Value = 1600.6
$ws.Range("J19").Value = 1774
$ws.Range("L19").Value = 1774
$ws.Range("N19").Value = -2124
# Row 53
$ws.Range("H53").Value = 326.84616
$ws.Range("I53").Value = 209.375
$ws.Range("K53").Value = 209.375
$ws.Range("M53").Value = 427.625
# Row 70
$ws.Range("H70").Value = 8432.25
$ws.Range("J70").Value = 7198
$ws.Range("L70").Value = 21594
$ws.Range("N70").Value = -22134
# Row 73
$ws.Range("H73").Value = 8432.25
$ws.Range("J73").Value = 7198
$ws.Range("L73").Value = 21594
$ws.Range("N73").Value = -23466
# Row 76
$ws.Range("H76").Value = 2003
$ws.Range("I76").Value = 2003
$ws.Range("K76").Value = 2003
$ws.Range("M76").Value = -1688
# Row 79
$ws.Range("H79").Value = 2003
$ws.Range("I79").Value = 2003
$ws.Range("K79").Value = 2003
$ws.Range("M79").Value = -911
# Row 80
$ws.Range("H80").Value = 798.8
$ws.Range("J80").Value = 699.5
$ws.Range("L80").Value = 2098.5
$ws.Range("N80").Value = -4094.5
# Row 83
$ws.Range("H83").Value = 798.8
$ws.Range("J83").Value = 699.5
$ws.Range("L83").Value = 6295.5
$ws.Range("N83").Value = -16279.5
# Row 113
$ws.Range("H113").Value = 3917.1667
$ws.Range("I113").Value = 3334.6667
$ws.Range("J113").Value = 4499.6665
$ws.Range("K113").Value = 3334.6667
$ws.Range("L113").Value = 4499.6665
$ws.Range("M113").Value = -80.66670000000022
$ws.Range("N113").Value = -11007.6665
# Row 132
$ws.Range("H132").Value = 1721.25
$ws.Range("I132").Value = 1759.3684
$ws.Range("K132").Value = 5278.1052
$ws.Range("M132").Value = -2748.1052
# Row 138
$ws.Range("H138").Value = 3436.7144
$ws.Range("J138").Value = 3726.7646
$ws.Range("L138").Value = 11180.2938
$ws.Range("N138").Value = -21460.2938

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 8
$ws.Range("H8").Value = 1146.2
$ws.Range("I8").Value = 1241.6666
$ws.Range("J8").Value = 1003
$ws.Range("K8").Value = 1241.6666
$ws.Range("L8").Value = 1003
$ws.Range("M8").Value = -1097.6666
$ws.Range("N8").Value = -1291
# Row 16
$ws.Range("H16").Value = 2007
$ws.Range("J16").Value = 2007
$ws.Range("L16").Value = 2007
$ws.Range("N16").Value = -2581
# Row 23
$ws.Range("H23").Value = 2000
$ws.Range("J23").Value = 2000
$ws.Range("L23").Value = 2000
$ws.Range("N23").Value = -2518
# Row 37
$ws.Range("H37").Value = 10499.5
$ws.Range("J37").Value = 10499.5
$ws.Range("L37").Value = 10499.5
$ws.Range("N37").Value = -11045.5
# Row 88
$ws.Range("H88").Value = 1609.3334
$ws.Range("I88").Value = 1000
$ws.Range("K88").Value = 1000
$ws.Range("M88").Value = -594
# Row 91
$ws.Range("H91").Value = 1609.3334
$ws.Range("I91").Value = 1000
$ws.Range("K91").Value = 1000
$ws.Range("M91").Value = 404
# Row 97
$ws.Range("H97").Value = 879.3333
$ws.Range("I97").Value = 892.46155
$ws.Range("J97").Value = 794
$ws.Range("K97").Value = 892.46155
$ws.Range("L97").Value = 794
$ws.Range("M97").Value = -396.46155
$ws.Range("N97").Value = -1786

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 1800.3334
$ws.Range("I22").Value = 5000
$ws.Range("J22").Value = 200.5
$ws.Range("K22").Value = 5000
$ws.Range("L22").Value = 200.5
$ws.Range("M22").Value = -4827
$ws.Range("N22").Value = -546.5
# Row 25
$ws.Range("H25").Value = 21000
$ws.Range("I25").Value = 10000
$ws.Range("K25").Value = 10000
$ws.Range("M25").Value = -9765
# Row 105
$ws.Range("H105").Value = 2733.75
$ws.Range("I105").Value = 2733.75
$ws.Range("K105").Value = 2733.75
$ws.Range("M105").Value = -986.75
# Row 107
$ws.Range("H107").Value = 3061.5
$ws.Range("I107").Value = 2830.5
$ws.Range("J107").Value = 3985.5
$ws.Range("K107").Value = 2830.5
$ws.Range("L107").Value = 3985.5
$ws.Range("M107").Value = -910.5
$ws.Range("N107").Value = -7825.5
# Row 125
$ws.Range("H125").Value = 83950
$ws.Range("J125").Value = 83950
$ws.Range("L125").Value = 83950
$ws.Range("N125").Value = -93790
# Row 134
$ws.Range("H134").Value = 3496.5833
$ws.Range("I134").Value = 3496.5833
$ws.Range("K134").Value = 10489.7499
$ws.Range("M134").Value = -7954.749899999999

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 5546.5557
$ws.Range("J16").Value = 5816.3335
$ws.Range("L16").Value = 5816.3335
$ws.Range("N16").Value = -6390.3335
# Row 106
$ws.Range("H106").Value = 65000
$ws.Range("J106").Value = 65000
$ws.Range("L106").Value = 65000
$ws.Range("N106").Value = -67524
# Row 113
$ws.Range("H113").Value = 5546.5557
$ws.Range("J113").Value = 5816.3335
$ws.Range("L113").Value = 5816.3335
$ws.Range("N113").Value = -10156.3335
# Row 114
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
# Row 132
$ws.Range("H132").Value = 5569.1763
$ws.Range("I132").Value = 5673.5
$ws.Range("K132").Value = 17020.5
$ws.Range("M132").Value = -14490.5
# Row 134
$ws.Range("H134").Value = 1391.5
$ws.Range("I134").Value = 1465.1666
$ws.Range("J134").Value = 949.5
$ws.Range("K134").Value = 4395.4998
$ws.Range("L134").Value = 2848.5
$ws.Range("M134").Value = -1860.4998
$ws.Range("N134").Value = -7918.5
# Row 135
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 32
$ws.Range("H32").Value = 326.33334
$ws.Range("J32").Value = 589
$ws.Range("L32").Value = 1767
$ws.Range("N32").Value = -2333
# Row 68
$ws.Range("H68").Value = 3829.9722
$ws.Range("J68").Value = 3905.1143
$ws.Range("L68").Value = 11715.3429
$ws.Range("N68").Value = -13337.3429
# Row 71
$ws.Range("H71").Value = 3829.9722
$ws.Range("J71").Value = 3905.1143
$ws.Range("L71").Value = 35146.0287
$ws.Range("N71").Value = -43258.0287
# Row 139
$ws.Range("H139").Value = 2677.6428
$ws.Range("I139").Value = 2432.3333
$ws.Range("K139").Value = 7296.999899999999
$ws.Range("M139").Value = -2156.999899999999

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 3863.3333
$ws.Range("I80").Value = 3456
$ws.Range("K80").Value = 3456
$ws.Range("M80").Value = -2458
# Row 83
$ws.Range("H83").Value = 3863.3333
$ws.Range("I83").Value = 3456
$ws.Range("K83").Value = 17280
$ws.Range("M83").Value = -12288

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 122
$ws.Range("H122").Value = 6483.0293
$ws.Range("I122").Value = 5029.8667
$ws.Range("K122").Value = 15089.6001
$ws.Range("M122").Value = -12639.6001

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 45
$ws.Range("H45").Value = 15721.5
$ws.Range("J45").Value = 15721.5
$ws.Range("L45").Value = 15721.5
$ws.Range("N45").Value = -16703.5
# Row 62
$ws.Range("H62").Value = 13423.75
$ws.Range("I62").Value = 14481.667
$ws.Range("J62").Value = 10250
$ws.Range("K62").Value = 14481.667
$ws.Range("L62").Value = 10250
$ws.Range("M62").Value = -13857.667
$ws.Range("N62").Value = -11498
# Row 65
$ws.Range("H65").Value = 13423.75
$ws.Range("I65").Value = 14481.667
$ws.Range("J65").Value = 10250
$ws.Range("K65").Value = 72408.33499999999
$ws.Range("L65").Value = 10250
$ws.Range("M65").Value = -69288.33499999999
$ws.Range("N65").Value = -57490
# Row 122
$ws.Range("H122").Value = 2553.6765
$ws.Range("I122").Value = 1777.2632
$ws.Range("K122").Value = 5331.7896
$ws.Range("M122").Value = -2881.7896
